$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.019.96'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.641.48'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('E4').Value = '  -0.49%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.71'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5091'
$ws.Range('E6').Value = '  +1.02%  '
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2567'
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06365'
$ws.Range('E9').Value = '  -0.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.59'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07780'
$ws.Range('E11').Value = '  -0.36%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.286'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.635.09'
$ws.Range('E13').Value = '  -0.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5445'
$ws.Range('E14').Value = '  +0.20%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0₅7744'
$ws.Range('E15').Value = '  -1.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.15'
$ws.Range('E16').Value = '  -1.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.042.20'
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('E18').Value = '  -0.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '196.76'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.423'
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.932'
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.036'
$ws.Range('E22').Value = '  +0.89%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.003'
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.873'
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.57'
$ws.Range('E25').Value = '  +0.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1192'
$ws.Range('E26').Value = '  +4.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.831'
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.53'
$ws.Range('E28').Value = '  -1.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.235'
$ws.Range('E29').Value = '  -0.65%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04865'
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.259'
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.170'
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.525'
$ws.Range('E33').Value = '  -0.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.361'
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.8976'
$ws.Range('E35').Value = '  +0.76%  '
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.579'
$ws.Range('E36').Value = '  -1.26%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.141.07'
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5470'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01563'
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.001'
$ws.Range('E40').Value = '  -0.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.520'
$ws.Range('E41').Value = '  -1.88%  '
$ws.Range('E42').Value = '  +8.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8094'
$ws.Range('E43').Value = '  -0.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.17'
$ws.Range('E44').Value = '  -0.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.400'
$ws.Range('E45').Value = '  -5.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.779.28'
$ws.Range('E46').Value = '  +0.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4527'
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.93'
$ws.Range('E48').Value = '  -0.82%  '
$ws.Range('B49').Value = 'Frax'
$ws.Range('C49').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.0000'
$ws.Range('E49').Value = '  -1.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05056'
$ws.Range('E50').Value = '  -0.67%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.001'
$ws.Range('E51').Value = '  -0.60%  '
